$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# The existing "Threshold" column (E) becomes "oldThreshold"; a new
# "Threshold" column is added in G (the new measurements), and a new
# "no object" column is added in J.
$ws.Range("E1").Value = "oldThreshold"

$ws.Range("G1").Value = "Threshold"
$ws.Range("G1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("J1").Value = "no object"
$ws.Range("J1").HorizontalAlignment = -4108  # xlCenter

# --- New "Threshold" (G) and "no object" (J) measurement columns ---
$ws.Range("G2").Value = 16.9
$ws.Range("J2").Value = 13.4

$ws.Range("G3").Value = 17.6
$ws.Range("J3").Value = 13.4

$ws.Range("G4").Value = 43.8
$ws.Range("J4").Value = 40.1

$ws.Range("G5").Value = 45.2
$ws.Range("J5").Value = 40.1

$ws.Range("G6").Value = 52.8
$ws.Range("J6").Value = 48.1

$ws.Range("G7").Value = 52.3
$ws.Range("J7").Value = 48.1

$ws.Range("G8").Value = 22.8
$ws.Range("J8").Value = 19.400000000000002

$ws.Range("G9").Value = 24.2
$ws.Range("J9").Value = 19.400000000000002

$ws.Range("G10").Value = 49.6
$ws.Range("J10").Value = 45.1

$ws.Range("G11").Value = 53.1
$ws.Range("J11").Value = 45.1

$ws.Range("G12").Value = 53.6
$ws.Range("J12").Value = 52.5

$ws.Range("G13").Value = 55.5
$ws.Range("J13").Value = 52.5

$ws.Range("G2:G13").HorizontalAlignment = -4108  # xlCenter

# --- New column of differences: Threshold - "no object" for each pair ---
$ws.Range("G14").Formula = "=G2-J2"
$ws.Range("G15").Formula = "=G3-J3"
$ws.Range("G16").Formula = "=G4-J4"
$ws.Range("G17").Formula = "=G5-J5"
$ws.Range("G18").Formula = "=G6-J6"
$ws.Range("G19").Formula = "=G7-J7"
$ws.Range("G20").Formula = "=G8-J8"
$ws.Range("G21").Formula = "=G9-J9"
$ws.Range("G22").Formula = "=G10-J10"
$ws.Range("G23").Formula = "=G11-J11"
$ws.Range("G24").Formula = "=G12-J12"
$ws.Range("G25").Formula = "=G13-J13"

$ws.Range("G14:G25").HorizontalAlignment = -4108  # xlCenter

# --- Clear stray empty formatted cells left over from earlier edits ---
$ws.Range("N14:Q14").Clear()
$ws.Range("Q15").Clear()
$ws.Range("J16:Q16").Clear()
$ws.Range("P17:Q17").Clear()
$ws.Range("P18:Q18").Clear()
$ws.Range("J19:Q19").Clear()
$ws.Range("J20:Q20").Clear()
$ws.Range("J21:K21").Clear()

# --- View/selection tidy-up ---
$ws.Range("J14").Select()
